$d = $word.ActiveDocument

# Update the date heading in the first paragraph
$d.Paragraphs(1).Range.Text = "2025-01-03 Friday"

# Update each table cell value by (row, column) position directly,
# so that duplicate expression text across cells is handled unambiguously.
$t = $d.Tables(1)

$t.Rows(1).Cells(1).Range.Text = "22+55="
$t.Rows(1).Cells(2).Range.Text = "6+3="
$t.Rows(1).Cells(3).Range.Text = "79-34="
$t.Rows(1).Cells(4).Range.Text = "30+58="
$t.Rows(1).Cells(5).Range.Text = "16+66="
$t.Rows(2).Cells(1).Range.Text = "68-14="
$t.Rows(2).Cells(2).Range.Text = "84-28="
$t.Rows(2).Cells(3).Range.Text = "80-40="
$t.Rows(2).Cells(4).Range.Text = "95-88="
$t.Rows(2).Cells(5).Range.Text = "63-20="
$t.Rows(3).Cells(1).Range.Text = "86-62="
$t.Rows(3).Cells(2).Range.Text = "83+9="
$t.Rows(3).Cells(3).Range.Text = "48-47="
$t.Rows(3).Cells(4).Range.Text = "56-23="
$t.Rows(3).Cells(5).Range.Text = "27+9="
$t.Rows(4).Cells(1).Range.Text = "35-17="
$t.Rows(4).Cells(2).Range.Text = "22+27="
$t.Rows(4).Cells(3).Range.Text = "98-71="
$t.Rows(4).Cells(4).Range.Text = "20+4="
$t.Rows(4).Cells(5).Range.Text = "72-11="
$t.Rows(5).Cells(1).Range.Text = "45-26="
$t.Rows(5).Cells(2).Range.Text = "92-57="
$t.Rows(5).Cells(3).Range.Text = "12+57="
$t.Rows(5).Cells(4).Range.Text = "66-62="
$t.Rows(5).Cells(5).Range.Text = "7+66="
$t.Rows(6).Cells(1).Range.Text = "19+75="
$t.Rows(6).Cells(2).Range.Text = "51+18="
$t.Rows(6).Cells(3).Range.Text = "30+55="
$t.Rows(6).Cells(4).Range.Text = "77-3="
$t.Rows(6).Cells(5).Range.Text = "82+11="
$t.Rows(7).Cells(1).Range.Text = "3+69="
$t.Rows(7).Cells(2).Range.Text = "6+63="
$t.Rows(7).Cells(3).Range.Text = "86-67="
$t.Rows(7).Cells(4).Range.Text = "81-44="
$t.Rows(7).Cells(5).Range.Text = "69-19="
$t.Rows(8).Cells(1).Range.Text = "59-54="
$t.Rows(8).Cells(2).Range.Text = "16+49="
$t.Rows(8).Cells(3).Range.Text = "29+14="
$t.Rows(8).Cells(4).Range.Text = "17+17="
$t.Rows(8).Cells(5).Range.Text = "76-63="
$t.Rows(9).Cells(1).Range.Text = "21+10="
$t.Rows(9).Cells(2).Range.Text = "27-15="
$t.Rows(9).Cells(3).Range.Text = "52+21="
$t.Rows(9).Cells(4).Range.Text = "28-2="
$t.Rows(9).Cells(5).Range.Text = "1+94="
$t.Rows(10).Cells(1).Range.Text = "60-59="
$t.Rows(10).Cells(2).Range.Text = "49-33="
$t.Rows(10).Cells(3).Range.Text = "9+0="
$t.Rows(10).Cells(4).Range.Text = "20+54="
$t.Rows(10).Cells(5).Range.Text = "34+26="
$t.Rows(11).Cells(1).Range.Text = "94+5="
$t.Rows(11).Cells(2).Range.Text = "77-3="
$t.Rows(11).Cells(3).Range.Text = "4+39="
$t.Rows(11).Cells(4).Range.Text = "13+37="
$t.Rows(11).Cells(5).Range.Text = "78+10="
$t.Rows(12).Cells(1).Range.Text = "60-49="
$t.Rows(12).Cells(2).Range.Text = "47-43="
$t.Rows(12).Cells(3).Range.Text = "26+55="
$t.Rows(12).Cells(4).Range.Text = "85-4="
$t.Rows(12).Cells(5).Range.Text = "0+61="
$t.Rows(13).Cells(1).Range.Text = "91-71="
$t.Rows(13).Cells(2).Range.Text = "97-20="
$t.Rows(13).Cells(3).Range.Text = "66-41="
$t.Rows(13).Cells(4).Range.Text = "38-29="
$t.Rows(13).Cells(5).Range.Text = "82-17="
$t.Rows(14).Cells(1).Range.Text = "85-7="
$t.Rows(14).Cells(2).Range.Text = "49+34="
$t.Rows(14).Cells(3).Range.Text = "40-15="
$t.Rows(14).Cells(4).Range.Text = "94-44="
$t.Rows(14).Cells(5).Range.Text = "20+29="
$t.Rows(15).Cells(1).Range.Text = "77-45="
$t.Rows(15).Cells(2).Range.Text = "49-6="
$t.Rows(15).Cells(3).Range.Text = "81+8="
$t.Rows(15).Cells(4).Range.Text = "25-3="
$t.Rows(15).Cells(5).Range.Text = "45+51="
$t.Rows(16).Cells(1).Range.Text = "78-41="
$t.Rows(16).Cells(2).Range.Text = "96-90="
$t.Rows(16).Cells(3).Range.Text = "38-3="
$t.Rows(16).Cells(4).Range.Text = "23+0="
$t.Rows(16).Cells(5).Range.Text = "89-73="
$t.Rows(17).Cells(1).Range.Text = "0+26="
$t.Rows(17).Cells(2).Range.Text = "85-31="
$t.Rows(17).Cells(3).Range.Text = "67-16="
$t.Rows(17).Cells(4).Range.Text = "61-48="
$t.Rows(17).Cells(5).Range.Text = "32+49="
$t.Rows(18).Cells(1).Range.Text = "80-47="
$t.Rows(18).Cells(2).Range.Text = "12+62="
$t.Rows(18).Cells(3).Range.Text = "72-67="
$t.Rows(18).Cells(4).Range.Text = "16+72="
$t.Rows(18).Cells(5).Range.Text = "55+40="
$t.Rows(19).Cells(1).Range.Text = "56+0="
$t.Rows(19).Cells(2).Range.Text = "73-10="
$t.Rows(19).Cells(3).Range.Text = "86-77="
$t.Rows(19).Cells(4).Range.Text = "7+10="
$t.Rows(19).Cells(5).Range.Text = "9+5="
$t.Rows(20).Cells(1).Range.Text = "51-48="
$t.Rows(20).Cells(2).Range.Text = "46+45="
$t.Rows(20).Cells(3).Range.Text = "80-19="
$t.Rows(20).Cells(4).Range.Text = "88-14="
$t.Rows(20).Cells(5).Range.Text = "95-80="
